# Apply the edit described by the diff:
# 1. Rename worksheet "Checklist" -> "Session"
# 2. Delete the first data row (Student ID 201007) shifting rows up
# 3. Change the "Type" column value from "Selection" to "Scan" for remaining rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Session"

# Delete row 2 (the 201007 log entry), which shifts subsequent rows up by one
$ws.Rows.Item(2).Delete()

# Update the Type column (E) from "Selection" to "Scan" for all remaining data rows (2-7)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 5).Value = "Scan"
}
